$d = $word.ActiveDocument
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function SetParaXML($paraIndex, $innerXml) {
    $p = $d.Paragraphs($paraIndex)
    $xml = "<w:p $wns>$innerXml</w:p>"
    $p.Range.InsertXML($xml) | Out-Null
}

$listBulletPPr = '<w:pPr><w:pStyle w:val="ListBullet"/><w:spacing w:line="240" w:lineRule="auto"/><w:ind w:left="720"/></w:pPr>'

# 1. Main title heading (Heading1 style)
$inner1 = '<w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:r><w:t>Play Mystery Express Free: Exciting Slot Game Review</w:t></w:r>'
SetParaXML 1 $inner1

# 2. "What we like" bullet list
$inner33 = $listBulletPPr + '<w:r/><w:r><w:t>Traditional structure with 30 pay lines</w:t></w:r>'
SetParaXML 33 $inner33

$inner34 = $listBulletPPr + '<w:r/><w:r><w:t>Interesting murder mystery storyline</w:t></w:r>'
SetParaXML 34 $inner34

$inner35 = $listBulletPPr + '<w:r/><w:r><w:t>Wide betting range catering to all players</w:t></w:r>'
SetParaXML 35 $inner35

$inner36 = $listBulletPPr + '<w:r/><w:r><w:t>Unique bonus features with exclusive rewards</w:t></w:r>'
SetParaXML 36 $inner36

# 3. "What we don't like" bullet list
$inner38 = $listBulletPPr + '<w:r/><w:r><w:t>Murder mystery is a side plot rather than a major focus</w:t></w:r>'
SetParaXML 38 $inner38

$inner39 = $listBulletPPr + '<w:r/><w:r><w:t>Limited number of pay lines</w:t></w:r>'
SetParaXML 39 $inner39

# 4. Bold meta title (repeated headline at bottom)
$inner40 = '<w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Mystery Express Free: Exciting Slot Game Review</w:t></w:r>'
SetParaXML 40 $inner40

# 5. Italic meta description
$inner41 = '<w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>Play Mystery Express for free and explore an adventurous Victorian London train ride. Review of gameplay, betting range, volatility, symbols, and bonus features.</w:t></w:r>'
SetParaXML 41 $inner41
